$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")
$summary = $wb.Worksheets.Item("Daily Summary")

# Insert a new row at row 2, shifting the existing orders down by one.
$ws.Rows.Item(2).Insert()

# Populate the new order row (row 2) with the new order's data.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "2026-01-13 18:48"
$ws.Range("C2").Value = "Sagar Borse"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "7588930329"
$ws.Range("E2").Value = "Test2,"
$ws.Range("F2").Value = "Kite Haldi Kunku Set x10"
$ws.Range("G2").Value = 300
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"

# Update the Daily Summary sheet totals to reflect the new order.
$summary.Range("B2").Value = 7
$summary.Range("E2").Value = 325
$summary.Range("G2").Value = 325
